$d = $word.ActiveDocument

$pairs = @(
    @("698×8=", "762×4="),
    @("925×4=", "638×5="),
    @("299×5=", "598×2="),
    @("912×9=", "447×7="),
    @("781×7=", "974×8="),
    @("221×6=", "569×7="),
    @("673×6=", "712×5="),
    @("485×9=", "681×4="),
    @("546×2=", "652×7="),
    @("519×5=", "170×5="),
    @("756×7=", "563×2="),
    @("772×2=", "507×4="),
    @("143×5=", "562×2="),
    @("980×5=", "650×2="),
    @("406×7=", "659×9="),
    @("791×5=", "707×3="),
    @("610×6=", "518×3="),
    @("879×8=", "490×8="),
    @("841×2=", "470×3="),
    @("663×9=", "300×5="),
    @("194×6=", "744×8="),
    @("440×6=", "753×3="),
    @("217×8=", "168×8="),
    @("144×9=", "659×7="),
    @("807×2=", "125×4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
